# Auto-generated edit script applying diff changes to cryptos worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns retain text formatting so
# values like "304.80" or "0.35%" are stored as literal strings, matching
# the source data (inline/shared strings), not auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "304.80"
$ws.Range("E2").Value = "0.35%"

# Row 3
$ws.Range("D3").Value = "35.90"
$ws.Range("E3").Value = "-3.42%"

# Row 4
$ws.Range("D4").Value = "5.068"
$ws.Range("E4").Value = "1.10%"

# Row 5
$ws.Range("D5").Value = "0.07906"
$ws.Range("E5").Value = "0.56%"

# Row 6
$ws.Range("D6").Value = "2.119"
$ws.Range("E6").Value = "-4.02%"

# Row 7
$ws.Range("D7").Value = "7.956"
$ws.Range("E7").Value = "-0.51%"

# Row 8
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "0.9235"
$ws.Range("E8").Value = "0.23%"

# Row 9
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "0.09761"
$ws.Range("E9").Value = "2.17%"

# Row 10
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1845"
$ws.Range("E10").Value = "-2.08%"

# Row 11
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "0.08667"
$ws.Range("E11").Value = "0.95%"

# Row 12
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "0.03568"
$ws.Range("E12").Value = "-1.03%"

# Row 13
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "0.09937"
$ws.Range("E13").Value = "-0.44%"

# Row 14
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "0.001444"
$ws.Range("E14").Value = "-1.60%"

# Row 15
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "0.005723"
$ws.Range("E15").Value = "1.00%"

# Row 16
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "3.461"
$ws.Range("E16").Value = "-0.08%"

# Row 17
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "4.131"
$ws.Range("E17").Value = "2.87%"

# Row 18
$ws.Range("E18").Value = "22.28%"

# Row 19
$ws.Range("D19").Value = "0.3381"
$ws.Range("E19").Value = "-1.12%"

# Row 20
$ws.Range("E20").Value = "2.40%"

# Row 21
$ws.Range("D21").Value = "5.169"
$ws.Range("E21").Value = "8.60%"

# Row 23
$ws.Range("D23").Value = "0.04563"
$ws.Range("E23").Value = "-0.61%"

# Row 24
$ws.Range("D24").Value = "0.001235"
$ws.Range("E24").Value = "0.43%"

# Row 25
$ws.Range("D25").Value = "0.004883"

# Row 26
$ws.Range("E26").Value = "-6.79%"

# Row 27
$ws.Range("D27").Value = "0.0004756"
$ws.Range("E27").Value = "0.11%"

# Row 39
$ws.Range("D39").Value = "0.01856"
$ws.Range("E39").Value = "1.13%"

# Row 40
$ws.Range("D40").Value = "0.04720"
$ws.Range("E40").Value = "-0.18%"

# Row 41
$ws.Range("D41").Value = "0.007921"
$ws.Range("E41").Value = "-2.52%"

# Row 42
$ws.Range("D42").Value = "0.1397"
$ws.Range("E42").Value = "0.18%"

# Row 43
$ws.Range("D43").Value = "0.007583"
$ws.Range("E43").Value = "0.51%"

# Row 44
$ws.Range("D44").Value = "0.002194"
$ws.Range("E44").Value = "-0.58%"

# Row 45
$ws.Range("D45").Value = "0.01127"
$ws.Range("E45").Value = "7.66%"

# Row 46
$ws.Range("D46").Value = "0.00006282"
$ws.Range("E46").Value = "1.86%"

# Row 47
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").Value = "0.28%"

# Row 48
$ws.Range("E48").Value = "0.25%"

# Row 49
$ws.Range("D49").Value = "49.94"
$ws.Range("E49").Value = "71.92%"

# Row 50
$ws.Range("D50").Value = "0.001903"
$ws.Range("E50").Value = "-29.27%"

# Row 51
$ws.Range("D51").Value = "0.00002103"
$ws.Range("E51").Value = "0.28%"
